$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 444.5
$ws.Range("J58").Value = 703.6
$ws.Range("L58").Value = 2110.8
$ws.Range("N58").Value = -2410.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 5478.7144
$ws.Range("I70").Value = 8251
$ws.Range("J70").Value = 4369.8
$ws.Range("K70").Value = 24753
$ws.Range("L70").Value = 13109.4
$ws.Range("M70").Value = -24483
$ws.Range("N70").Value = -13649.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 5478.7144
$ws.Range("I73").Value = 8251
$ws.Range("J73").Value = 4369.8
$ws.Range("K73").Value = 24753
$ws.Range("L73").Value = 13109.4
$ws.Range("M73").Value = -23817
$ws.Range("N73").Value = -14981.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 2185.375
$ws.Range("I74").Value = 1783.2858
$ws.Range("K74").Value = 1783.2858
$ws.Range("M74").Value = -847.2858000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 2185.375
$ws.Range("I77").Value = 1783.2858
$ws.Range("K77").Value = 8916.429
$ws.Range("M77").Value = -4236.429

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 2999
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").Value = $null

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3596.25
$ws.Range("I32").Value = 3596.25
$ws.Range("K32").Value = 3596.25
$ws.Range("M32").Value = -3309.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 6859
$ws.Range("I45").Value = 2045
$ws.Range("J45").Value = 8062.5
$ws.Range("K45").Value = 2045
$ws.Range("L45").Value = 8062.5
$ws.Range("M45").Value = -1668
$ws.Range("N45").Value = -8816.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 4042.75
$ws.Range("J63").Value = 2500
$ws.Range("L63").Value = 2500
$ws.Range("N63").Value = -3872

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 4042.75
$ws.Range("J66").Value = 2500
$ws.Range("L66").Value = 12500
$ws.Range("N66").Value = -19364

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H95").Value = 16933.166
$ws.Range("J95").Value = 16933.166
$ws.Range("L95").Value = 16933.166
$ws.Range("N95").Value = -22425.166

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").Value = $null

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1241.2
$ws.Range("I86").Value = 999.6667
$ws.Range("J86").Value = 1603.5
$ws.Range("K86").Value = 999.6667
$ws.Range("L86").Value = 1603.5
$ws.Range("M86").Value = 123.3333
$ws.Range("N86").Value = -3849.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1241.2
$ws.Range("I89").Value = 999.6667
$ws.Range("J89").Value = 1603.5
$ws.Range("K89").Value = 4998.3335
$ws.Range("L89").Value = 8017.5
$ws.Range("M89").Value = 617.6665000000003
$ws.Range("N89").Value = -19249.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 900
$ws.Range("I94").Value = 900
$ws.Range("K94").Value = 900
$ws.Range("M94").Value = -449

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1700
$ws.Range("I105").Value = 1700
$ws.Range("K105").Value = 1700
$ws.Range("M105").Value = 47

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 70.71429000000001
$ws.Range("I19").Value = 70.71429000000001
$ws.Range("K19").Value = 70.71429000000001
$ws.Range("M19").Value = 99.28570999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H24").Value = 70.71429000000001
$ws.Range("I24").Value = 70.71429000000001
$ws.Range("K24").Value = 70.71429000000001
$ws.Range("M24").Value = 99.28570999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 15000
$ws.Range("I25").Value = 15000
$ws.Range("K25").Value = 15000
$ws.Range("M25").Value = -14826

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8608.666999999999
$ws.Range("I31").Value = 5691.231
$ws.Range("K31").Value = 5691.231
$ws.Range("M31").Value = -5396.231

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 8608.666999999999
$ws.Range("I34").Value = 5691.231
$ws.Range("K34").Value = 5691.231
$ws.Range("M34").Value = -5489.231

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 10326.25
$ws.Range("I58").Value = 5402.5
$ws.Range("K58").Value = 5402.5
$ws.Range("M58").Value = -5199.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 14998
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").Value = $null

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H103").Value = 0
$ws.Range("I103").Value = 0
$ws.Range("K103").Value = 0
$ws.Range("M103").Value = $null

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 14998
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").Value = $null

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 6436
$ws.Range("I132").Value = 5372
$ws.Range("K132").Value = 16116
$ws.Range("M132").Value = -13586

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 8190.625
$ws.Range("I134").Value = 2506
$ws.Range("K134").Value = 7518
$ws.Range("M134").Value = -4983

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H135").Value = 74988.664
$ws.Range("I135").Value = 74987
$ws.Range("J135").Value = 74989.5
$ws.Range("K135").Value = 74987
$ws.Range("L135").Value = 74989.5
$ws.Range("M135").Value = -69917
$ws.Range("N135").Value = -85129.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 10326.25
$ws.Range("I136").Value = 5402.5
$ws.Range("K136").Value = 16207.5
$ws.Range("M136").Value = -13657.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 3690.9092
$ws.Range("J11").Value = 3690.9092
$ws.Range("L11").Value = 11072.7276
$ws.Range("N11").Value = -11352.7276

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H119").Value = 0
$ws.Range("I119").Value = 0
$ws.Range("K119").Value = 0
$ws.Range("M119").Value = $null

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = $null
$ws.Range("N47").Value = 0

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H104").Value = 30671
$ws.Range("J104").Value = 30671
$ws.Range("L104").Value = 30671
$ws.Range("N104").Value = -37659

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H105").Value = 19557
$ws.Range("J105").Value = 19557
$ws.Range("L105").Value = 19557
$ws.Range("N105").Value = -26545

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2169.3333
$ws.Range("I16").Value = 2088
$ws.Range("J16").Value = 2210
$ws.Range("K16").Value = 2088
$ws.Range("L16").Value = 2210
$ws.Range("M16").Value = -1918
$ws.Range("N16").Value = -2550

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 2692.7778
$ws.Range("J55").Value = 2847
$ws.Range("L55").Value = 2847
$ws.Range("N55").Value = -3193

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2200
$ws.Range("I93").Value = 2200
$ws.Range("K93").Value = 2200
$ws.Range("M93").Value = -952

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = $null
$ws.Range("M100").Value = $null
$ws.Range("N100").Value = 0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = $null
$ws.Range("N108").Value = 0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H139").Value = 131571.33
$ws.Range("I139").Value = 98000
$ws.Range("J139").Value = 148357
$ws.Range("K139").Value = 98000
$ws.Range("L139").Value = 148357
$ws.Range("M139").Value = -92860
$ws.Range("N139").Value = -158637

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 40
$ws.Range("I13").Value = 40
$ws.Range("K13").Value = 40
$ws.Range("M13").Value = 100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 31400
$ws.Range("J54").Value = 31400
$ws.Range("L54").Value = 31400
$ws.Range("N54").Value = -32440

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 12498
$ws.Range("I62").Value = 12498
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 12498
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = $null
$ws.Range("N62").Value = -11874

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 12498
$ws.Range("I65").Value = 12498
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 62490
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = $null
$ws.Range("N65").Value = -59370

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 956.2857
$ws.Range("I100").Value = 239.25
$ws.Range("K100").Value = 478.5
$ws.Range("M100").Value = 62.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H103").Value = 23534
$ws.Range("J103").Value = 23534
$ws.Range("L103").Value = 23534
$ws.Range("N103").Value = -25878

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 24068.8
$ws.Range("J104").Value = 25086
$ws.Range("L104").Value = 25086
$ws.Range("N104").Value = -32074

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H105").Value = 40000
$ws.Range("J105").Value = 40000
$ws.Range("L105").Value = 40000
$ws.Range("N105").Value = -46988

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1773.25
$ws.Range("I113").Value = 1156.6
$ws.Range("K113").Value = 3469.8
$ws.Range("M113").Value = -1299.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2872.75
$ws.Range("I122").Value = 3498.8333
$ws.Range("J122").Value = 994.5
$ws.Range("K122").Value = 10496.4999
$ws.Range("L122").Value = 2983.5
$ws.Range("M122").Value = -8046.499899999999
$ws.Range("N122").Value = -7883.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 13134.667
$ws.Range("I136").Value = 3404
$ws.Range("J136").Value = 18000
$ws.Range("K136").Value = 10212
$ws.Range("L136").Value = 54000
$ws.Range("M136").Value = -7662
$ws.Range("N136").Value = -59100
